# Feat/49/transfer crud: add a new "TRANSFER" worksheet between ACCOUNT and
# REGULAR_TRANSFER, and add a TRANSFER_ID column to TEMPORARY_TRANSFER.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new TRANSFER sheet right after ACCOUNT (i.e. before the
#    existing REGULAR_TRANSFER sheet), which becomes sheetId 6 / sheet3.xml.
# ------------------------------------------------------------------
$accountSheet = $wb.Worksheets.Item("ACCOUNT")
$transferSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $accountSheet)
$transferSheet.Name = "TRANSFER"

# Header row
$transferSheet.Range("A1").Value = "ID"
$transferSheet.Range("B1").Value = "TITLE"
$transferSheet.Range("B1").NumberFormat = "@"
$transferSheet.Range("C1").Value = "USER_ID"

# Data rows
$transferSheet.Range("A2").Value = 1
$transferSheet.Range("B2").Value = "2023/06"
$transferSheet.Range("B2").NumberFormat = "@"
$transferSheet.Range("C2").Value = 1

$transferSheet.Range("A3").Value = 2
$transferSheet.Range("B3").Value = "2023/07"
$transferSheet.Range("B3").NumberFormat = "@"
$transferSheet.Range("C3").Value = 1

$transferSheet.Range("A4").Value = 3
$transferSheet.Range("B4").Value = "2022/09"
$transferSheet.Range("B4").NumberFormat = "@"
$transferSheet.Range("C4").Value = 2

# ------------------------------------------------------------------
# 2. Add a TRANSFER_ID column (G) to TEMPORARY_TRANSFER.
# ------------------------------------------------------------------
$temporaryTransferSheet = $wb.Worksheets.Item("TEMPORARY_TRANSFER")
$temporaryTransferSheet.Range("G1").Value = "TRANSFER_ID"
$temporaryTransferSheet.Range("G2").Value = 1
$temporaryTransferSheet.Range("G3").Value = 2
$temporaryTransferSheet.Range("G4").Value = 3

# ------------------------------------------------------------------
# 3. Selections / active sheet, so the saved view matches: TRANSFER's
#    selection becomes A1:C4, and TEMPORARY_TRANSFER ends up the active
#    (selected) tab with H14 selected.
# ------------------------------------------------------------------
$transferSheet.Range("A1:C4").Select()
$temporaryTransferSheet.Range("H14").Select()
